$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the judge's non-availability date in D8 (2018-07-31 -> 2018-07-30)
$ws.Range("D8").Value = 43311

# Move/save the active cell selection to C15 (matches author's last selection on save)
[void]$ws.Range("C15").Select()
